# Refresh cryptocurrency price/volume data (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.359.37"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'1.848.84"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'240.44"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'0.6289"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.07612"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.2920"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'24.62"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'0.07746"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "'5.021"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'0.6810"
$ws.Range("D14").Value = "'0.00001050"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "'83.13"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'6.130"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'29.351.55"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'229.55"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'12.33"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "'0.9998"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'7.472"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'158.51"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").Value = "'0.1391"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'8.441"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +10.14%  "
$ws.Range("D28").Value = "'1.474"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "'0.05618"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").Value = "'4.111"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'4.055"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "'1.831"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'0.6986"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").Value = "'2.582"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.01812"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "'1.235.25"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D39").Value = "'6.416"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").Value = "'0.9015"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").Value = "'0.9992"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'101.56"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "'65.62"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "'7.169"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1158"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.029"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.685"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000114"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("D50").Value = "'0.05701"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'0.4629"
$ws.Range("E51").Value = "  -0.07%  "
